$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: the "License Information" heading paragraph becomes a plain
# (non-heading) paragraph whose bold run now reads
# "Aquifer Open Study Notes (Book Intros)".
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(4)
$frag1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
    '<w:pPr><w:bidi/></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="ar_AE" w:bidi="ar_AE"/></w:rPr></w:r>' +
    '<w:r><w:rPr><w:b/><w:lang w:val="ar_AE" w:bidi="ar_AE"/></w:rPr><w:t>Aquifer Open Study Notes (Book Intros)</w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'
$p1.Range.InsertXML($frag1)

# ---------------------------------------------------------------------------
# Edit 2: the attribution paragraph is rewritten to describe the new
# Aquifer Open Study Notes adaptation instead of the old Tyndale wording,
# and the two hyperlinks are removed entirely.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(5)
$frag2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
    '<w:pPr><w:bidi/></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="ar_AE" w:bidi="ar_AE"/></w:rPr></w:r>' +
    '<w:r><w:rPr><w:lang w:val="ar_AE" w:bidi="ar_AE"/></w:rPr><w:t xml:space="preserve">This work is an adaptation of </w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="ar_AE" w:bidi="ar_AE"/></w:rPr><w:t>Tyndale Open Study Notes</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="ar_AE" w:bidi="ar_AE"/></w:rPr><w:t xml:space="preserve"> © 2023 Tyndale House Publishers, licensed under the CC BY-SA 4.0 license. The adaptation, </w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="ar_AE" w:bidi="ar_AE"/></w:rPr><w:t>Aquifer Open Study Notes</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="ar_AE" w:bidi="ar_AE"/></w:rPr><w:t>, was created by Mission Mutual and is also licensed under CC BY-SA 4.0.</w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'
$p2.Range.InsertXML($frag2)

# ---------------------------------------------------------------------------
# Edit 3: the closing sentence now lists the languages the resource has
# been adapted into, instead of the old single-sentence PDF note. Use
# InsertXML (not Find/Replace) so the paragraph's other runs (including the
# two trailing empty ones) are preserved exactly.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(6)
$frag3 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
    '<w:pPr><w:bidi/></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="ar_AE" w:bidi="ar_AE"/></w:rPr></w:r>' +
    '<w:r><w:rPr><w:lang w:val="ar_AE" w:bidi="ar_AE"/></w:rPr><w:t>This resource has been adapted into multiple languages, including English, Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文).</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="ar_AE" w:bidi="ar_AE"/></w:rPr></w:r>' +
    '<w:r><w:rPr><w:lang w:val="ar_AE" w:bidi="ar_AE"/></w:rPr></w:r>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'
$p3.Range.InsertXML($frag3)
